$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows
$ws.Range("F6").Value = -7
$ws.Range("F7").Value = 5
$ws.Range("F10").Value = -8
$ws.Range("F13").Value = -8
$ws.Range("F20").Value = -9
$ws.Range("F21").Value = -3
$ws.Range("F24").Value = -3
$ws.Range("F26").Value = -2
$ws.Range("F28").Value = -5
